# Insert a new data record at row 3 (right after the header-adjacent first
# data row), pushing the existing rows 3-72 down to 4-73, then populate the
# newly inserted row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 3; Excel shifts rows 3:72 down to 4:73
# and copies the formatting (incl. the date style on column D) from the row
# above, same as native Excel UI behaviour.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new record.
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(3, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(3, 4).Value = 44496
$ws.Cells.Item(3, 5).Value = 15
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100102
$ws.Cells.Item(3, 8).Value = "Cítricos"
$ws.Cells.Item(3, 9).Value = 100102005
$ws.Cells.Item(3, 10).Value = "Naranja"
$ws.Cells.Item(3, 11).Value = "Lane Late"
$ws.Cells.Item(3, 12).Value = "Segunda"
$ws.Cells.Item(3, 13).Value = 250
$ws.Cells.Item(3, 14).Value = 650
$ws.Cells.Item(3, 15).Value = 700
$ws.Cells.Item(3, 16).Value = 675
$ws.Cells.Item(3, 17).Value = '$/kilo (en caja de 20 kilos)'
$ws.Cells.Item(3, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(3, 19).Value = 675
$ws.Cells.Item(3, 20).Value = 1
